# Rename the "How to implement basic daily attendance" Heading1 section to
# "How to implement basic financial transactions" - this means updating both
# the visible heading text AND the bookmark that marks that heading
# (the bookmark's w:name is used for in-document anchors / TOC links).

$d = $word.ActiveDocument

$oldText = "How to implement basic daily attendance"
$newText = "How to implement basic financial transactions"
$oldBookmarkName = "how-to-implement-basic-daily-attendance"
$newBookmarkName = "how-to-implement-basic-financial-transactions"

# Locate the Heading1 paragraph that holds the old title.
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text.TrimEnd([char]13, [char]7) -eq $oldText) {
        $target = $p
        break
    }
}

if ($target -eq $null) {
    throw "Could not find the '$oldText' heading paragraph"
}

$oldRange = $target.Range
$oldStart = $oldRange.Start
$oldEnd = $oldRange.End

# Insert a brand-new paragraph (same style as the one that follows it, i.e.
# Heading1) immediately before the old one, carrying the new title text.
$insertionPoint = $d.Range($oldStart, $oldStart)
$insertionPoint.InsertBefore($newText + [char]13)

# Re-create the bookmark (collapsed, right at the start of the heading, just
# like the original) under its new name on the freshly inserted paragraph.
$newBookmarkRange = $d.Range($oldStart, $oldStart)
$d.Bookmarks.Add($newBookmarkName, $newBookmarkRange)

# Remove the now-duplicate old heading paragraph (text + its bookmark) - its
# start/end shifted forward by the length of the text + paragraph mark we
# just inserted.
$shift = $newText.Length + 1
$oldParagraphRange = $d.Range($oldStart + $shift, $oldEnd + $shift)
$oldParagraphRange.Delete()

Write-Output "Renamed heading '$oldText' -> '$newText' (bookmark '$oldBookmarkName' -> '$newBookmarkName')"
